$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (values chosen so the stored/rounded width lands
# exactly on the target integer character width after Excel's internal
# pixel rounding)
$ws.Columns.Item(3).ColumnWidth = 40.166666666666664
$ws.Columns.Item(4).ColumnWidth = 37.166666666666664
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 20.166666666666668

$data = @(
    @("1331206", "https://aiesec.org/opportunity/global-talent/1331206", "Consulting Intern", "Peshawar, Pakistan", "No", "2 applicants", "9 - 12 Weeks", "iConsult"),
    @("1331205", "https://aiesec.org/opportunity/global-talent/1331205", "Finance Intern Controller", "Panamá, Provincia de Panamá, Panamá", "No", "1 applicant", "6 - 18 Months", "HILTI Panama"),
    @("1331181", "https://aiesec.org/opportunity/global-talent/1331181", "Marketing Executive", "Cyberjaya, Selangor, Malaysia", "No", "9 applicants", "6 - 18 Months", "IX Telecom Sdn Bhd"),
    @("1331055", "https://aiesec.org/opportunity/global-talent/1331055", "Repair technician", "Hammam Sousse, Tunisie", "No", "1 applicant", "3 - 6 Months", "MOBYSTORE"),
    @("1328766", "https://aiesec.org/opportunity/global-talent/1328766", "Accelerate Romania - Community Manager", "Cluj-Napoca, Romania", "No", "35 applicants", "9 - 12 Weeks", "ClujStartups"),
    @("1321294", "https://aiesec.org/opportunity/global-talent/1321294", "UI/UX & Web Development Intern", "Lahore, Punjab, Pakistan", "No", "44 applicants", "9 - 12 Weeks", "Devsinc."),
    @("1312369", "https://aiesec.org/opportunity/global-talent/1312369", "Software Engineer Intern", "Lahore, Punjab, Pakistan", "No", "8 applicants", "9 - 12 Weeks", "Devsinc"),
    @("1309734", "https://aiesec.org/opportunity/global-talent/1309734", "Marketing Intern", "Lahore, Punjab, Pakistan", "No", "12 applicants", "9 - 12 Weeks", "Devsinc.")
)

$rowIndex = 2
foreach ($rowData in $data) {
    $colIndex = 1
    foreach ($val in $rowData) {
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        if ($colIndex -eq 1) {
            # Opportunity ID looks numeric - force text storage like the source data
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
        $colIndex++
    }
    $rowIndex++
}
